# Mifos -> Finflux: 1st changes
# On the "Repayment schedule" sheet, a new (blank) column is inserted
# before column N ("Late"), pushing the existing N/O/P columns
# ("Late"/"heading"/"Outstanding") one position to the right (O/P/Q).
# The active sheet/tab also moves from "Transactions" to "Repayment schedule",
# with a new selection of S4 on that sheet.

$wb = $excel.ActiveWorkbook

$wsRepay = $wb.Worksheets.Item("Repayment schedule")

# Insert a blank column at position 14 (N); everything from N onward shifts right.
$wsRepay.Columns.Item(14).Insert() | Out-Null

# The newly inserted column picks up the width of the column immediately to
# its left (M, which holds the original "In Advance" column's width).
$wsRepay.Columns.Item(14).ColumnWidth = $wsRepay.Columns.Item(13).ColumnWidth

# Move the active sheet/selection: "Repayment schedule" becomes the
# selected/active tab (previously it was "Transactions"), with cell S4
# selected on it.
$wsRepay.Activate() | Out-Null
$wsRepay.Range("S4").Select() | Out-Null
